$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ψ (rad)"
$ws.Range("E1").Value = "ψ (deg)"
$ws.Range("B2").Value = 0
$ws.Range("D2").Value = -0.8106398939946861
$ws.Range("E2").Value = -46.44624463082796
$ws.Range("A3").Value = 0.0362158009982007
$ws.Range("B3").Value = 0.0909546412732028
$ws.Range("C3").Value = 2.404333723078434
$ws.Range("D3").Value = -0.8106398939946861
$ws.Range("E3").Value = -46.44624463082796
$ws.Range("A4").Value = 0.07243160199640139
$ws.Range("B4").Value = 0.1819092825464056
$ws.Range("C4").Value = 2.308667446156868
$ws.Range("D4").Value = -0.8094300737272562
$ws.Range("E4").Value = -46.37692703553484
$ws.Range("A5").Value = 0.1086474029946021
$ws.Range("B5").Value = 0.2729795962285458
$ws.Range("C5").Value = 2.213111277988574
$ws.Range("D5").Value = -0.8070100603747469
$ws.Range("E5").Value = -46.23827048407075
$ws.Range("A6").Value = 0.1448632039928028
$ws.Range("B6").Value = 0.3642808902129477
$ws.Range("C6").Value = 2.117775780790811
$ws.Range("D6").Value = -0.8033791078421669
$ws.Range("E6").Value = -46.03023222834158
$ws.Range("A7").Value = 0.1810790049910035
$ws.Range("B7").Value = 0.4559277402526472
$ws.Range("C7").Value = 2.022772421971391
$ws.Range("D7").Value = -0.798536095836598
$ws.Range("E7").Value = -45.7527480802913
$ws.Range("A8").Value = 0.2172948059892042
$ws.Range("B8").Value = 0.5480336161255041
$ws.Range("C8").Value = 1.928214022350496
$ws.Range("D8").Value = -0.7924795284832962
$ws.Range("E8").Value = -45.40573233261038
$ws.Range("A9").Value = 0.2535106069874049
$ws.Range("B9").Value = 0.6407104985041194
$ws.Range("C9").Value = 1.834215199050675
$ws.Range("D9").Value = -0.7852075324748449
$ws.Range("E9").Value = -44.98907765269014
$ws.Range("A10").Value = 0.2897264079856056
$ws.Range("B10").Value = 0.7340684834707161
$ws.Range("C10").Value = 1.74089280113813
$ws.Range("D10").Value = -0.7767178547484879
$ws.Range("E10").Value = -44.50265494954366
$ws.Range("A11").Value = 0.3259422089838063
$ws.Range("B11").Value = 0.8282153716522893
$ws.Range("C11").Value = 1.648366335978712
$ws.Range("D11").Value = -0.7670078596855029
$ws.Range("E11").Value = -43.94631321334176
$ws.Range("A12").Value = 0.362158009982007
$ws.Range("B12").Value = 0.9232562389979436
$ws.Range("C12").Value = 1.556758384122137
$ws.Range("D12").Value = -0.7560745258251982
$ws.Range("E12").Value = -43.31987932713882
$ws.Range("A13").Value = 0.3983738109802076
$ws.Range("B13").Value = 1.019292986280172
$ws.Range("C13").Value = 1.466195000348008
$ws.Range("D13").Value = -0.7439144420847872
$ws.Range("E13").Value = -42.62315785028762
$ws.Range("A14").Value = 0.4345896119784083
$ws.Range("B14").Value = 1.116423864477102
$ws.Range("C14").Value = 1.376806098297117
$ws.Range("D14").Value = -0.7305238034750351
$ws.Range("E14").Value = -41.85593077296389
$ws.Range("A15").Value = 0.4708054129766091
$ws.Range("B15").Value = 1.214742973285989
$ws.Range("C15").Value = 1.288725815871383
$ws.Range("D15").Value = -0.7158984063001664
$ws.Range("E15").Value = -41.01795724114136
$ws.Range("A16").Value = 0.5070212139748097
$ws.Range("B16").Value = 1.31433973013249
$ws.Range("C16").Value = 1.202092858315786
$ws.Range("D16").Value = -0.7000336428290562
$ws.Range("E16").Value = -40.10897325127343
$ws.Range("A17").Value = 0.5432370149730105
$ws.Range("B17").Value = 1.415298307178871
$ws.Range("C17").Value = 1.117050815596078
$ws.Range("D17").Value = -0.6829244954232082
$ws.Range("E17").Value = -39.12869131385114
$ws.Range("A18").Value = 0.5794528159712111
$ws.Range("B18").Value = 1.517697034001182
$ws.Range("C18").Value = 1.033748450357507
$ws.Range("D18").Value = -0.6645655301054324
$ws.Range("E18").Value = -38.07680008491553
$ws.Range("A19").Value = 0.6156686169694118
$ws.Range("B19").Value = 1.621607763804822
$ws.Range("C19").Value = 0.9523399523929471
$ws.Range("D19").Value = -0.6449508895514652
$ws.Range("E19").Value = -36.95296396450706
$ws.Range("A20").Value = 0.6518844179676125
$ws.Range("B20").Value = 1.727095201284628
$ws.Range("C20").Value = 0.8729851551648348
$ws.Range("D20").Value = -0.6240742854850243
$ws.Range("E20").Value = -35.75682266093435
$ws.Range("A21").Value = 0.6881002189658132
$ws.Range("B21").Value = 1.834216190514862
$ws.Range("C21").Value = 0.7958497095155953
$ws.Range("D21").Value = -0.6019289904549435
$ws.Range("E21").Value = -34.48799071963868
$ws.Range("A22").Value = 0.724316019964014
$ws.Range("B22").Value = 1.943018961581972
$ws.Range("C22").Value = 0.7211052092678258
$ws.Range("D22").Value = -0.5785078289710793
$ws.Range("E22").Value = -33.1460570153189
$ws.Range("A23").Value = 0.7605318209622146
$ws.Range("B23").Value = 2.053542335054975
$ws.Range("C23").Value = 0.6489292629607778
$ws.Range("D23").Value = -0.5538031679736147
$ws.Range("E23").Value = -31.73058420586272
$ws.Range("A24").Value = 0.7967476219604153
$ws.Range("B24").Value = 2.165814883831479
$ws.Range("C24").Value = 0.5795055054969369
$ws.Range("D24").Value = -0.5278069066081913
$ws.Range("E24").Value = -30.24110814650496
$ws.Range("A25").Value = 0.832963422958616
$ws.Range("B25").Value = 2.279854052408832
$ws.Range("C25").Value = 0.5130235429856163
$ws.Range("D25").Value = -0.5005104652769643
$ws.Range("E25").Value = -28.67713726249919
$ws.Range("A26").Value = 0.8691792239568167
$ws.Range("B26").Value = 2.395665234217337
$ws.Range("C26").Value = 0.4496788235743806
$ws.Range("D26").Value = -0.4719047739331893
$ws.Range("E26").Value = -27.03815187844697
$ws.Range("A27").Value = 0.9053950249550173
$ws.Range("B27").Value = 2.513240808323908
$ws.Range("C27").Value = 0.3896724265597211
$ws.Range("D27").Value = -0.441980259584287
$ws.Range("E27").Value = -25.3236035022762
$ws.Range("A28").Value = 0.9416108259532181
$ws.Range("B28").Value = 2.632559137578396
$ws.Range("C28").Value = 0.3332107615728607
$ws.Range("D28").Value = -0.4107268329654952
$ws.Range("E28").Value = -23.5329140616976
$ws.Range("A29").Value = 0.9778266269514188
$ws.Range("B29").Value = 2.753583531139724
$ws.Range("C29").Value = 0.2805051691534261
$ws.Range("D29").Value = -0.3781338743431645
$ws.Range("E29").Value = -21.66547509079353
$ws.Range("A30").Value = 1.014042427949619
$ws.Range("B30").Value = 2.876261175293962
$ws.Range("C30").Value = 0.2317714135631078
$ws.Range("D30").Value = -0.3441902184034877
$ws.Range("E30").Value = -19.72064686420588
$ws.Range("A31").Value = 1.05025822894782
$ws.Range("B31").Value = 3.000522037570412
$ws.Range("C31").Value = 0.1872290582652156
$ws.Range("D31").Value = -0.3088841381789377
$ws.Range("E31").Value = -17.69775747618887
$ws.Range("A32").Value = 1.086474029946021
$ws.Range("B32").Value = 3.126277750383729
$ws.Range("C32").Value = 0.1471007141181505
$ws.Range("D32").Value = -0.2722033279609016
$ws.Range("E32").Value = -15.59610186157506
$ws.Range("A33").Value = 1.122689830944221
$ws.Range("B33").Value = 3.253420481788657
$ws.Range("C33").Value = 0.1116111500173442
$ws.Range("D33").Value = -0.2341348851429262
$ws.Range("E33").Value = -13.41494075546995
$ws.Range("A34").Value = 1.158905631942422
$ws.Range("B34").Value = 3.381821802437416
$ws.Range("C34").Value = 0.08098625548976654
$ws.Range("D34").Value = -0.1946652909345802
$ws.Range("E34").Value = -11.15349958823773
$ws.Range("A35").Value = 1.195121432940623
$ws.Range("B35").Value = 3.511331559485492
$ws.Range("C35").Value = 0.05545184461896328
$ws.Range("D35").Value = -0.1537803898811876
$ws.Range("E35").Value = -8.810967312068364
$ws.Range("A36").Value = 1.231337233938824
$ws.Range("B36").Value = 3.641776770005932
$ws.Range("C36").Value = 0.03523229068103254
$ws.Range("D36").Value = -0.1114653681195327
$ws.Range("E36").Value = -6.3864951551213
$ws.Range("A37").Value = 1.267553034937024
$ws.Range("B37").Value = 3.772960548450054
$ws.Range("C37").Value = 0.02054898103046932
$ws.Range("D37").Value = -0.06770473029406145
$ws.Range("E37").Value = -3.87919529892125
$ws.Range("A38").Value = 1.303768835935225
$ws.Range("B38").Value = 3.904661084836499
$ws.Range("C38").Value = 0.01161858212037695
$ws.Range("D38").Value = -0.0224822750520533
$ws.Range("E38").Value = -1.288139474334917
$ws.Range("A39").Value = 1.339984636933426
$ws.Range("B39").Value = 4.036630692661015
$ws.Range("C39").Value = 0.008651105108825401
$ws.Range("D39").Value = 0.02421893097033651
$ws.Range("E39").Value = 1.387642528918962
$ws.Range("A40").Value = 1.376200437931626
$ws.Range("B40").Value = 4.16859494799322
$ws.Range("C40").Value = 0.01184776333072056
$ws.Range("D40").Value = 0.07241658076640689
$ws.Range("E40").Value = 4.149164444683367
$ws.Range("A41").Value = 1.412416238929827
$ws.Range("B41").Value = 4.300251943856724
$ws.Range("C41").Value = 0.02139861404509238
$ws.Range("D41").Value = 0.1221291562846347
$ws.Range("E41").Value = 6.997485210603202
$ws.Range("A42").Value = 1.448632039928028
$ws.Range("B42").Value = 4.431271686763377
$ws.Range("C42").Value = 0.0374799783480709
$ws.Range("D42").Value = 0.1733759595732904
$ws.Range("E42").Value = 9.933710752580323
